$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data (row 14)
$ws.Range("A14").Value = "09.07.2019 - " + [char]10 + "11.07.2019"
$ws.Range("B14").Value = "Researching on handling url parameters and query strings with get requests. Researching on how to do a bit more complex searches with SQLAlchemy. Working towards better queries and more robust database without PickleTypes which can't be used for filtering. Added the cards' colors as a separate table which can now be used as a filter."
$ws.Range("C14").Value = 12

# Apply wrap-text styling (same style used elsewhere, e.g. B13) to the new cells
$ws.Range("A14:B14").WrapText = $true

$ws.Rows.Item(14).RowHeight = 90

# Update the view/selection to match the target state
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("B17").Select() | Out-Null

$wb.Save()
